$wb = $excel.ActiveWorkbook

# The two sheets that need the new "hydrogen if" subscript rows added
$sheetNames = @("IFSR-chemicals", "ISFR-ironsteel")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 12: green hydrogen if
    $ws.Cells.Item(12, 1).Value = "green hydrogen if"
    $ws.Cells.Item(12, 2).Value = 1

    # Row 13: low carbon hydrogen if
    $ws.Cells.Item(13, 1).Value = "low carbon hydrogen if"
    $ws.Cells.Item(13, 2).Value = 1
}

# Update selections to mirror where the author ended up navigating
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("H19").Select()

$wsChem = $wb.Worksheets.Item("IFSR-chemicals")
$wsChem.Range("A14").Select()

$wsSteel = $wb.Worksheets.Item("ISFR-ironsteel")
$wsSteel.Activate()
$wsSteel.Range("A14").Select()
